# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $savedStyle
}

Set-TextValue "D2" "62.053.18"
$ws.Range("E2").Value = "  -2.07%  "
Set-TextValue "D3" "3.421.05"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "578.29"
$ws.Range("E5").Value = "  -0.59%  "
Set-TextValue "D6" "152.28"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.22%  "
Set-TextValue "D9" "8.06"
$ws.Range("E9").Value = "  +4.95%  "
$ws.Range("E10").Value = "  -0.70%  "
Set-TextValue "D11" "0.418"
$ws.Range("E11").Value = "  +3.29%  "
Set-TextValue "D12" "4.008.13"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -3.10%  "
Set-TextValue "D15" "3.417.37"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("E16").Value = "  -0.48%  "
Set-TextValue "D17" "62.107.79"
$ws.Range("E17").Value = "  -2.06%  "
Set-TextValue "D18" "6.51"
$ws.Range("E18").Value = "  +2.16%  "
Set-TextValue "D19" "14.58"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  -4.18%  "
Set-TextValue "D21" "382.75"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  +0.87%  "
Set-TextValue "D23" "75.17"
$ws.Range("E24").Value = "  +0.01%  "
Set-TextValue "D25" "3.563.10"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("E26").Value = "  -3.30%  "
Set-TextValue "D27" "0.180"
$ws.Range("E27").Value = "  -0.48%  "
Set-TextValue "D28" "7.69"
$ws.Range("E28").Value = "  +0.73%  "
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -1.12%  "
Set-TextValue "D31" "7.90"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  -1.95%  "
Set-TextValue "D34" "23.23"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D38" "31.10"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D39" "168.52"
$ws.Range("E39").Value = "  -0.45%  "
Set-TextValue "D40" "3.453.86"
$ws.Range("E40").Value = "  -1.74%  "
Set-TextValue "D41" "0.0785"
$ws.Range("E41").Value = "  +2.73%  "
Set-TextValue "D42" "42.74"
$ws.Range("E42").Value = "  +0.83%  "
Set-TextValue "D43" "0.779"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("E46").Value = "  -2.78%  "
Set-TextValue "D47" "2.542.24"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  -5.08%  "
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  +0.01%  "
